# Updated symbol list with refreshed price / 1h-volume-change figures
# for the crypto tracker sheet (cryptos.xlsx), Coin list rows 2-50.
#
# The Price (D) and Volume(1h) (E) columns are stored as plain text in
# this sheet (not numbers / percentages), so we force a Text number
# format while writing the new value and then restore the cell's
# original style -- this keeps Excel from "helpfully" re-typing
# "300.98" as a float or "0.13%" as a percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "300.98"
    "E2" = "0.13%"
    "D3" = "31.43"
    "E3" = "-0.29%"
    "D4" = "5.085"
    "E4" = "-1.05%"
    "D5" = "0.07863"
    "E5" = "-1.22%"
    "D6" = "2.347"
    "E6" = "-1.78%"
    "D7" = "7.811"
    "E7" = "-1.58%"
    "D8" = "3.841"
    "E8" = "-0.33%"
    "E9" = "0.94%"
    "D10" = "0.1751"
    "E10" = "1.68%"
    "D11" = "0.07587"
    "E11" = "3.04%"
    "D12" = "0.09216"
    "E12" = "14.42%"
    "D13" = "0.02999"
    "E13" = "-2.95%"
    "D14" = "0.1004"
    "E14" = "0.76%"
    "D15" = "0.001518"
    "E15" = "-0.17%"
    "D16" = "0.005845"
    "E16" = "-3.94%"
    "D17" = "3.473"
    "E17" = "-0.66%"
    "D18" = "2.247"
    "E18" = "0.26%"
    "D20" = "0.1290"
    "E20" = "-4.71%"
    "D21" = "4.112"
    "E21" = "-11.28%"
    "D22" = "0.1710"
    "E22" = "6.32%"
    "D23" = "0.04609"
    "E23" = "-0.27%"
    "E24" = "-1.07%"
    "D25" = "0.004465"
    "E25" = "0.12%"
    "D26" = "0.0001248"
    "E26" = "4.47%"
    "D27" = "0.0003387"
    "E27" = "-1.88%"
    "D39" = "0.01743"
    "E39" = "-5.89%"
    "D40" = "0.04673"
    "E40" = "3.18%"
    "D41" = "0.007062"
    "E41" = "-3.28%"
    "E42" = "1.21%"
    "D43" = "0.002186"
    "E43" = "0.60%"
    "D44" = "0.009760"
    "E44" = "-8.24%"
    "D45" = "0.00006286"
    "E45" = "-1.92%"
    "D46" = "0.00000000750"
    "E46" = "-0.43%"
    "E47" = "19.48%"
    "D48" = "0.7459"
    "E48" = "-9.10%"
    "D49" = "0.00002099"
    "E49" = "-0.43%"
    "D50" = "0.0001999"
    "E50" = "-0.43%"
}

foreach ($address in $updates.Keys) {
    $newValue = $updates[$address]
    $cell = $ws.Range($address)

    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = $originalStyle
}
